# Update the "Förändrad" (Changed) date column (C) for all data rows
# from serial date 46060 (2026-02-07) to 46061 (2026-02-08).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46060) {
        $cell.Value2 = 46061
    }
}
